$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.393000000000001
$ws.Range("D6").Value = -8.264000000000001
$ws.Range("D7").Value = -7.523999999999999
$ws.Range("D16").Value = -7.896999999999998
$ws.Range("D20").Value = -8.071999999999999
